# refactor the format of the xlsx
# Swap column B <-> C and column D <-> E (values, formatting and widths),
# matching how Excel performs a "cut column, insert before" reorder.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap columns B and C (Class 1 name / Class 1 code)
$ws.Columns.Item(3).Cut()
$ws.Columns.Item(2).Insert()

# Swap columns D and E (Class 2 name / Class 2 code)
$ws.Columns.Item(5).Cut()
$ws.Columns.Item(4).Insert()

# Update the active selection to the full column D, as left after the reorder
[void]$ws.Columns.Item(4).Select()
